# TestSuiteRestServices.xlsx edit
# 1. Prefix the TESTCASE values in column C with "service".
# 2. Move the active-cell selection from D6 to C6.
# 3. Nudge the A:E column widths down slightly (matches the resaved layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "serviceupload"
$ws.Range("C3").Value = "serviceinsightingest"
$ws.Range("C4").Value = "serviceextractDocumentMetadata"
$ws.Range("C5").Value = "serviceconvertDocument"
$ws.Range("C6").Value = "serviceclassifyDocument"
$ws.Range("C7").Value = "serviceextractDocumentElements"
$ws.Range("C8").Value = "serviceextractDocumentText"
$ws.Range("C9").Value = "serviceupload-serviceinsightingest"

$ws.Columns("A:A").ColumnWidth = 7.4173
$ws.Columns("B:B").ColumnWidth = 15.5894
$ws.Columns("C:C").ColumnWidth = 36.2524
$ws.Columns("D:D").ColumnWidth = 17.7576
$ws.Columns("E:E").ColumnWidth = 12.2596

$ws.Range("C6").Select() | Out-Null
